# Remove the 4th (last, blank) slide from the deck.
# This corresponds to sldId 259 / r:id rId5 / ppt/slides/slide4.xml
# being dropped from the presentation entirely.
$p = $ppt.ActivePresentation
$p.Slides.Item(4).Delete()
